$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Existing row 6 (FAN_CURRENT) DLC changes from 2 to 1
$ws.Range("E6").Value = 1

# New test rows for every parameter type - text cells entered in the same
# order the original author typed them so shared-string ids line up
$ws.Range("B7").Value = "U8_TESTER"
$ws.Range("B8").Value = "U16_TESTER"
$ws.Range("D9").Value = "U32"
$ws.Range("D10").Value = "U64"
$ws.Range("D11").Value = "S8"
$ws.Range("D12").Value = "S16"
$ws.Range("D13").Value = "S32"
$ws.Range("D14").Value = "S64"
$ws.Range("D15").Value = "FLOAT"
$ws.Range("B9").Value = "U32_TESTER"
$ws.Range("B11").Value = "S8_TESTER"
$ws.Range("B12").Value = "S16_TESTER"
$ws.Range("B13").Value = "S32_TESTER"
$ws.Range("B14").Value = "S64_TESTER"
$ws.Range("B15").Value = "FLOAT_TESTER"
$ws.Range("B10").Value = "U64_TESTER"

# remaining text cells that reuse already-existing shared strings
$ws.Range("D7").Value = "U8"
$ws.Range("D8").Value = "U16"

# numeric columns (ID / DLC)
$ws.Range("C7").Value = 4
$ws.Range("C8").Value = 5
$ws.Range("C9").Value = 6
$ws.Range("C10").Value = 7
$ws.Range("C11").Value = 8
$ws.Range("C12").Value = 9
$ws.Range("C13").Value = 10
$ws.Range("C14").Value = 11
$ws.Range("C15").Value = 12

$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 2
$ws.Range("E9").Value = 4
$ws.Range("E10").Value = 8
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 4
$ws.Range("E14").Value = 8
$ws.Range("E15").Value = 4

# wrap-text formatting carried down onto a couple of ID cells and a blank row
$ws.Range("C8").WrapText = $true
$ws.Range("C12").WrapText = $true
$ws.Range("C16").WrapText = $true

# update selection to reflect where the edit left off
$ws.Range("C18").Select()
